# Updates crypto price/volume data as scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.200.19"
$ws.Range("E2").Value = "  +6.18%  "
$ws.Range("D3").Value = "3.662.91"
$ws.Range("E3").Value = "  +17.66%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.50"
$ws.Range("E5").Value = "  +3.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.08"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("D7").Value = "3.658.97"
$ws.Range("E7").Value = "  +17.61%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  +4.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +7.78%  "
$ws.Range("E11").Value = "  +3.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.500"
$ws.Range("E12").Value = "  +5.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.87"
$ws.Range("E13").Value = "  +12.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000256"
$ws.Range("E14").Value = "  +5.72%  "
$ws.Range("D15").Value = "4.264.79"
$ws.Range("E15").Value = "  +17.28%  "
$ws.Range("D16").Value = "71.188.65"
$ws.Range("E16").Value = "  +6.25%  "
$ws.Range("D17").Value = "3.656.81"
$ws.Range("E17").Value = "  +17.16%  "
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.53"
$ws.Range("E19").Value = "  +7.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.06"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "522.69"
$ws.Range("E21").Value = "  +7.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.21"
$ws.Range("E22").Value = "  +18.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.747"
$ws.Range("E23").Value = "  +7.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.18"
$ws.Range("E24").Value = "  +5.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +10.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.61"
$ws.Range("E26").Value = "  +6.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.01"
$ws.Range("E27").Value = "  +6.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +10.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.14"
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.80"
$ws.Range("E31").Value = "  +7.31%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.89"
$ws.Range("E32").Value = "  +13.15%  "
$ws.Range("E33").Value = "  +17.20%  "
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.14"
$ws.Range("E36").Value = "  +9.01%  "
$ws.Range("E37").Value = "  +7.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.346"
$ws.Range("E38").Value = "  +11.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.18"
$ws.Range("E39").Value = "  +9.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.17"
$ws.Range("E40").Value = "  +3.77%  "
$ws.Range("E41").Value = "  +5.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.04"
$ws.Range("E42").Value = "  -8.13%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.85"
$ws.Range("E43").Value = "  +6.07%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.133.53"
$ws.Range("E44").Value = "  +12.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "416.03"
$ws.Range("E45").Value = "  +11.29%  "
$ws.Range("E46").Value = "  +4.45%  "
$ws.Range("E47").Value = "  +6.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.49"
$ws.Range("E48").Value = "  +13.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "139.71"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  +10.92%  "
